$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill in the new values for the "Hardwareansteuerung" sub-section (rows 22-29)
# Leading apostrophe -> stored as text with quotePrefix (matches how Excel
# marks text that starts with a "-" so it isn't mistaken for a formula/number)
$ws.Range("C22").Value = "'---------"
$ws.Range("C23").Value = "??"
$ws.Range("C24").Value = "'---------"
$ws.Range("C27").Value = 5
$ws.Range("C28").Value = 0

# Recalculate formulas (SUM(C22:C29) in D30, SUM(C2:C35) in C37)
$excel.Calculate()

# Update the view state to match where the user scrolled/selected afterwards
$ws.Activate()
$ws.Application.ActiveWindow.Zoom = 85
$ws.Application.ActiveWindow.ScrollRow = 13
[void]$ws.Range("C28").Select()
